# NIT-9015069037 Estado de Cuenta — update per commit:
# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# Changes:
#  1. Insert a new worker row (EDITH SIERRA ARENAS) at the top of the
#     data table (row 16), pushing the existing rows down.
#  2. Re-order the KAREN AURIMAR VALERA MONTERO rows so they come after
#     JOSE ALBERTO LAMO CASTAÑEDA instead of interleaved with JOSE MANUEL.
#  3. Update the "Valor Mora" total and the "Cant. Trabajadores" /
#     "Cant. Periodos" counters.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- 1. Insert a fresh row for the new worker right above the first
#        existing data row (row 16), shifting everything below it
#        (including the footer rows) down by one. ---
$ws.Rows("16:16").Insert()

# The freshly inserted row has no explicit formatting yet — copy the
# look of the (now shifted-down) data row directly beneath it so it
# matches the rest of the table exactly.
$ws.Range("B17:J17").Copy()
$ws.Range("B16:J16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# --- 2. Write the final contents of every data row (16-21). ---
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1143339752"
$ws.Range("D16").Value = "EDITH SIERRA ARENAS"
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 34164
$ws.Range("G16").Value = 1423500

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1235045585"
$ws.Range("D17").Value = "JOSE MANUEL CASTRO MONCARIS"
$ws.Range("E17").Value = "2211"
$ws.Range("F17").Value = 40000
$ws.Range("G17").Value = 1160000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1235045585"
$ws.Range("D18").Value = "JOSE MANUEL CASTRO MONCARIS"
$ws.Range("E18").Value = "2210"
$ws.Range("F18").Value = 40000
$ws.Range("G18").Value = 1160000

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1143263187"
$ws.Range("D19").Value = "JOSE ALBERTO LAMO CASTAÑEDA"
$ws.Range("E19").Value = "2305"
$ws.Range("F19").Value = 3712
$ws.Range("G19").Value = 1392000

$ws.Range("B20").Value = "PPT"
$ws.Range("C20").Value = "3959222"
$ws.Range("D20").Value = "KAREN AURIMAR VALERA MONTERO"
$ws.Range("E20").Value = "2211"
$ws.Range("F20").Value = 40000
$ws.Range("G20").Value = 1400000

$ws.Range("B21").Value = "PPT"
$ws.Range("C21").Value = "3959222"
$ws.Range("D21").Value = "KAREN AURIMAR VALERA MONTERO"
$ws.Range("E21").Value = "2210"
$ws.Range("F21").Value = 40000
$ws.Range("G21").Value = 1400000

# --- 3. Update the summary fields above the table. ---
$ws.Range("E11").Value = 197876      # VALOR MORA total
$ws.Range("C13").Value = 4           # Cant. Trabajadores
$ws.Range("F13").Value = 4           # Cant. Periodos
